# Daily update at 8 AM UTC
# Appends the new day's results (2025-06-29 / serial 45837) as row 98,
# and restores the previous last row (97) to the regular date format
# since it is no longer the final row in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 97 was previously the last row and used the "last row" date format
# (YYYY-MM-DD). Now that a new row follows it, it reverts to the normal
# date+time format used by all the other data rows.
$ws.Range("A97").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 98: Day / Chase / Bryce / Zach
$ws.Range("A98").Value = 45837
$ws.Range("A98").NumberFormat = "YYYY-MM-DD"
$ws.Range("B98").Value = 418
$ws.Range("C98").Value = 413
$ws.Range("D98").Value = 427
